$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains exact text representation (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.883.27"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "1.860.36"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "246.57"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").Value = "0.6371"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "0.3017"
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("D9").Value = "0.07492"
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("D10").Value = "24.55"
$ws.Range("E10").Value = "  +7.59%  "
$ws.Range("D11").Value = "0.07679"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "1.852.85"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").Value = "5.060"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "0.6904"
$ws.Range("E14").Value = "  +4.76%  "
$ws.Range("D15").Value = "84.46"
$ws.Range("E15").Value = "  +3.74%  "
$ws.Range("D16").Value = "0.000009435"
$ws.Range("E16").Value = "  +5.34%  "
$ws.Range("D17").Value = "6.105"
$ws.Range("E17").Value = "  +4.78%  "
$ws.Range("D18").Value = "29.835.68"
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("D19").Value = "2.117.74"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "240.37"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "12.69"
$ws.Range("E21").Value = "  +2.25%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "7.373"
$ws.Range("E23").Value = "  +3.78%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "159.00"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").Value = "0.1426"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").Value = "8.572"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("D28").Value = "18.00"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("D29").Value = "1.508"
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").Value = "0.06046"
$ws.Range("E30").Value = "  +8.72%  "
$ws.Range("D31").Value = "1.268"
$ws.Range("E31").Value = "  +5.08%  "
$ws.Range("D32").Value = "4.152"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "4.150"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").Value = "1.867"
$ws.Range("D35").Value = "1.163"
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("D36").Value = "0.7311"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "2.619"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "2.864"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").Value = "1.223.46"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D40").Value = "0.01791"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").Value = "6.331"
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "0.9226"
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.017.60"
$ws.Range("E44").Value = "  +2.86%  "
$ws.Range("D45").Value = "102.39"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").Value = "66.61"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").Value = "0.5086"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").Value = "9.320"
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("D50").Value = "0.4098"
$ws.Range("E50").Value = "  +3.01%  "
$ws.Range("D51").Value = "0.1143"
$ws.Range("E51").Value = "  +3.23%  "
